# Split the bold "JPMorgan Chase & Co." run into two runs: "JPMorgan" and
# "Chase" (dropping " & Co." entirely), per the target diff.

$d = $word.ActiveDocument

# Shrink the original bold run down to just "JPMorgan" (drops " Chase & Co.").
$rng = $d.Content
$rng.Find.Execute("JPMorgan Chase & Co.", $true, $false, $false, $false, $false, $true, 1, $false, "JPMorgan", 2)

# Re-locate "JPMorgan" and collapse to a point right after it.
$afterJPMorgan = $d.Content
$afterJPMorgan.Find.Execute("JPMorgan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterJPMorgan.Collapse(0)

# Insert "Chase" immediately after - same (inherited) bold formatting.
$afterJPMorgan.InsertAfter("Chase")

# Toggling Bold off/on over just the inserted text forces Word to keep it as
# its own run (identical <w:b/><w:bCs/> formatting) instead of silently
# re-merging it with the preceding "JPMorgan" run.
$chaseRng = $d.Range($afterJPMorgan.Start, $afterJPMorgan.Start + 5)
$chaseRng.Bold = 0
$chaseRng.Bold = 1
